$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextValue($addr, $val) {
    $c = $ws.Range($addr)
    $origStyle = $c.Style
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.Style = $origStyle
}

$ws.Range("D2").Value = '62.762.91'
$ws.Range("E2").Value = '  -0.37%  '

$ws.Range("D3").Value = '2.464.15'
$ws.Range("E3").Value = '  +0.38%  '

Set-TextValue "D5" '574.56'
$ws.Range("E5").Value = '  -0.68%  '

Set-TextValue "D6" '147.92'
$ws.Range("E6").Value = '  +1.01%  '

$ws.Range("E7").Value = '  +0.01%  '

Set-TextValue "D8" '0.533'
$ws.Range("E8").Value = '  -1.27%  '

$ws.Range("D9").Value = '2.465.98'
$ws.Range("E9").Value = '  +0.52%  '

Set-TextValue "D10" '0.111'
$ws.Range("E10").Value = '  -0.39%  '

$ws.Range("E11").Value = '  -0.49%  '

Set-TextValue "D12" '5.28'
$ws.Range("E12").Value = '  -0.35%  '

Set-TextValue "D13" '0.356'
$ws.Range("E13").Value = '  +0.21%  '

Set-TextValue "D14" '29.18'
$ws.Range("E14").Value = '  +2.36%  '

$ws.Range("E15").Value = '  -1.72%  '

$ws.Range("D16").Value = '2.915.65'
$ws.Range("E16").Value = '  +0.59%  '

$ws.Range("D17").Value = '62.691.17'
$ws.Range("E17").Value = '  -0.34%  '

$ws.Range("D18").Value = '2.465.92'
$ws.Range("E18").Value = '  +0.42%  '

Set-TextValue "D19" '7.93'
$ws.Range("E19").Value = '  -0.81%  '

Set-TextValue "D20" '10.99'
$ws.Range("E20").Value = '  -1.01%  '

Set-TextValue "D21" '326.95'
$ws.Range("E21").Value = '  -1.39%  '

Set-TextValue "D23" '2.20'
$ws.Range("E23").Value = '  +5.75%  '

$ws.Range("E24").Value = '  -0.02%  '

$ws.Range("E25").Value = '  +17.32%  '

Set-TextValue "D26" '65.52'
$ws.Range("E26").Value = '  -1.35%  '

Set-TextValue "D27" '639.58'
$ws.Range("E27").Value = '  -1.71%  '

$ws.Range("D28").Value = '0.0₃0981'
$ws.Range("E28").Value = '  -2.23%  '

Set-TextValue "D30" '0.992'
$ws.Range("E30").Value = '  -21.25%  '

$ws.Range("E31").Value = '  -0.97%  '

Set-TextValue "D32" '7.94'
$ws.Range("E32").Value = '  -3.35%  '

Set-TextValue "D33" '1.83'
$ws.Range("E33").Value = '  -2.44%  '

$ws.Range("E34").Value = '  -2.57%  '

Set-TextValue "D35" '0.999'
$ws.Range("E35").Value = '  +0.04%  '

$ws.Range("E36").Value = '  +3.21%  '

Set-TextValue "D37" '4.74'
$ws.Range("E37").Value = '  -0.88%  '

$ws.Range("B38").Value = 'PolygonEcosystemToken'
$ws.Range("C38").Value = 'https://coinranking.com/coin/iDZ0tG-wI+polygonecosystemtoken-pol'
Set-TextValue "D38" '0.368'
$ws.Range("E38").Value = '  -1.72%  '

$ws.Range("B39").Value = 'Monero'
$ws.Range("C39").Value = 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr'
Set-TextValue "D39" '151.36'
$ws.Range("E39").Value = '  -0.60%  '

$ws.Range("B40").Value = 'EthereumClassic'
$ws.Range("C40").Value = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
Set-TextValue "D40" '18.68'
$ws.Range("E40").Value = '  -0.80%  '

$ws.Range("B41").Value = 'dogwifhat'
$ws.Range("C41").Value = 'https://coinranking.com/coin/sZUrmToWF+dogwifhat-wif'
Set-TextValue "D41" '2.79'
$ws.Range("E41").Value = '  +1.47%  '

$ws.Range("B42").Value = 'RenderToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/vfo5XYwcV+rendertoken-render'
Set-TextValue "D42" '5.36'
$ws.Range("E42").Value = '  -3.41%  '

Set-TextValue "D43" '1.73'
$ws.Range("E43").Value = '  -2.35%  '

$ws.Range("E44").Value = '  -0.02%  '

$ws.Range("D45").Value = '0.0₆0304'
$ws.Range("E45").Value = '  -28.83%  '

Set-TextValue "D46" '153.14'
$ws.Range("E46").Value = '  +4.55%  '

$ws.Range("E47").Value = '  +1.71%  '

Set-TextValue "D48" '3.58'
$ws.Range("E48").Value = '  -1.60%  '

Set-TextValue "D49" '20.42'
$ws.Range("E49").Value = '  -1.63%  '

$ws.Range("E50").Value = '  +0.34%  '

$ws.Range("E51").Value = '  -1.48%  '
